$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$hours = @("05","07","08","09","10","11","12","13","14","15","16","17","18","19","20","21","22")

for ($i = 0; $i -lt $hours.Length; $i++) {
    $row = 3 + $i
    $ws.Range("B$row").Value = "$($hours[$i]):00"
}
